# Edit script for Enunt-laborator02.docx
# Applies the renumbering / rewording of lab exercises 3-9 plus adds
# "(8_findNum.c) " and "(9_findNumAuto.c) " labels to two trailing items.
#
# NOTE: In this runtime, Find.Execute(..., Replace:=2/wdReplaceOne) behaves
# like "replace all" across the searched range, not "replace first". To
# stay safe we always locate the (first) match with a non-replacing Find
# and then set .Text directly on the located Range, which only touches
# that one occurrence.

$d = $word.ActiveDocument

function Find-First($searchText) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        Write-Output "WARNING: not found: $searchText"
        return $null
    }
    return $rng
}

function Replace-First($searchText, $newText) {
    $rng = Find-First $searchText
    if ($rng -ne $null) {
        $rng.Text = $newText
    }
    return $rng
}

# ---------------------------------------------------------------------
# Item 3 (send100.c) -> Item 2 (send.c), plus reworded body text.
# ---------------------------------------------------------------------
Replace-First "3_" "2_" | Out-Null
Replace-First "send100.c" "send.c" | Out-Null
Replace-First ") Modificați programul anterior în așa fel încât în loc de un element să fie transmis un vector de " ") Modificați programul anterior adăugând transmisia unui vector de " | Out-Null
Replace-First "de elemente o dată, printr-un singur apel. Aveți grijă să inițializați vectorul doar pe procesul 0." "de elemente. Se va executa întreaga transmisie printr-un singur apel. Aveți grijă să inițializați vectorul doar pe procesul 0." | Out-Null

Write-Output "Item 3->2 done"

# ---------------------------------------------------------------------
# Item 4 (broadcast.c) -> Item 3 (broadcast.c) - number only.
# ---------------------------------------------------------------------
Replace-First "4_" "3_" | Out-Null

# ---------------------------------------------------------------------
# Item 5 (broadcast100.c) -> Item 3 (broadcast.c), plus reworded body.
# ---------------------------------------------------------------------
Replace-First "5_" "3_" | Out-Null
Replace-First "broadcast100.c" "broadcast.c" | Out-Null
Replace-First ") Modificați programul anterior în așa fel încât în loc de un element să fie transmis un vector de 100 de elemente o dată, printr-un singur apel." ") Modificați programul anterior în așa fel încât să adăugați transmisia unui vector de 100 de element. Se va executa întreaga transmisie printr-un singur apel." | Out-Null

Write-Output "Item 4->3 and 5->3 done"

# ---------------------------------------------------------------------
# Items 6,7,8,9 -> 4,5,6,7 - number only.
# ---------------------------------------------------------------------
Replace-First "6_" "4_" | Out-Null
Replace-First "7_" "5_" | Out-Null
Replace-First "8_" "6_" | Out-Null
Replace-First "9_" "7_" | Out-Null

Write-Output "Items 6-9 renumbered"
